$d = $word.ActiveDocument

$replacements = @(
    @("22÷9=2, 4", "91÷8=11, 3"),
    @("84÷5=16, 4", "71÷7=10, 1"),
    @("13÷6=2, 1", "99÷3=33, 0"),
    @("32÷2=16, 0", "14÷8=1, 6"),
    @("73÷9=8, 1", "84÷6=14, 0"),
    @("48÷3=16, 0", "34÷5=6, 4"),
    @("18÷7=2, 4", "47÷6=7, 5"),
    @("59÷7=8, 3", "74÷7=10, 4"),
    @("31÷6=5, 1", "86÷3=28, 2"),
    @("85÷4=21, 1", "22÷2=11, 0"),
    @("18÷8=2, 2", "54÷5=10, 4"),
    @("54÷7=7, 5", "56÷7=8, 0"),
    @("98÷9=10, 8", "93÷5=18, 3"),
    @("18÷2=9, 0", "41÷6=6, 5"),
    @("77÷4=19, 1", "35÷8=4, 3"),
    @("44÷9=4, 8", "26÷4=6, 2"),
    @("57÷9=6, 3", "42÷9=4, 6"),
    @("97÷2=48, 1", "30÷5=6, 0"),
    @("23÷4=5, 3", "30÷3=10, 0"),
    @("47÷4=11, 3", "55÷5=11, 0"),
    @("77÷5=15, 2", "80÷8=10, 0"),
    @("84÷9=9, 3", "75÷3=25, 0"),
    @("89÷8=11, 1", "38÷6=6, 2"),
    @("79÷7=11, 2", "94÷7=13, 3"),
    @("32÷8=4, 0", "33÷5=6, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
